# Generate Report for Handback
# Updates the handback status report with freshly-generated timestamps and
# a status change for the 5b82cd21-... / cc745ee5-... rows.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the
# 5b82cd21-... (row 3) and cc745ee5-... (row 5) entries.
$wsOverview.Range("G3").Value = "2016-08-19 04:13:34"
$wsOverview.Range("G5").Value = "2016-08-19 04:13:34"

# zh-cn sheet: Status (E), Correspond Handoff Datetime (H) and
# Correspond Handback DateTime (K) for the same two rows.
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"

$wsZhCn.Range("H3").Value = "2016-08-19 04:13:26"
$wsZhCn.Range("H5").Value = "2016-08-19 04:13:26"

$wsZhCn.Range("K3").Value = "2016-08-19 04:13:54"
$wsZhCn.Range("K5").Value = "2016-08-19 04:13:54"

# de-de sheet: Status (E), Correspond Handoff Datetime (H) and
# Correspond Handback DateTime (K) for the same two rows.
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"

$wsDeDe.Range("H3").Value = "2016-08-19 04:13:34"
$wsDeDe.Range("H5").Value = "2016-08-19 04:13:34"

$wsDeDe.Range("K3").Value = "2016-08-19 04:14:03"
$wsDeDe.Range("K5").Value = "2016-08-19 04:14:03"
